$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FBS")

# --- Remove the two games that dropped off the slate ---
# Row 8 = "Penn State @ Oregon", Row 5 = "Iowa State @ Arizona State"
# Delete the lower row first so the upper row index stays valid.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()

# After the two deletions the remaining games shift up to rows 2-7:
#   2 Tulane @ Army
#   3 Western Kentucky @ Jacksonville State
#   4 UNLV @ Boise State
#   5 Ohio @ Miami (OH)
#   6 Marshall @ Louisiana
#   7 Clemson @ SMU

# --- Refresh odds / wind values that changed between snapshots ---

# Row 3: Western Kentucky @ Jacksonville State
$ws.Range("AB3").Value = -4.5
$ws.Range("AF3").Value = 1

# Row 4: UNLV @ Boise State
$ws.Range("Q4").Value = "SW"
$ws.Range("Y4").Value = 57.5
$ws.Range("Z4").Value = -105
$ws.Range("AE4").Value = -0.0170940170940171

# Row 7: Clemson @ SMU
$ws.Range("Y7").Value = 55.5
$ws.Range("Z7").Value = -115
$ws.Range("AB7").Value = -2
$ws.Range("AE7").Value = 0.01834862385321101
$ws.Range("AF7").Value = -0.5

# --- Stamp every remaining data row with the new run's timestamp ---
$timestamp = "2024-12-04T16:21:35.504286"
for ($r = 2; $r -le 7; $r++) {
    $ws.Range("AK" + $r).Value = $timestamp
}
